$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (TC 003): userName/password updated, Expected updated.
# Entering these first matches the original author's edit order, which
# determines shared-string insertion order.
$ws.Range("B4").Value = "'hemanthgeneraluser@unilogcorp.com"
$ws.Range("C4").Value = "'hemanth123"

# Row 2 (TC 001): userName/password updated.
$ws.Range("B2").Value = "hemanth.bs123@unilogcorp.com"
$ws.Range("C2").Value = "'hemanth123"

# Row 4 Expected updated last.
$ws.Range("D4").Value = "Purchase Agent"

# Row 3 (TC 002): userName/password/Expected text is unchanged, but the
# mailto hyperlink on B3 is removed along with its "Hyperlink" style.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B3").Style = "Normal"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:hemanth.bs123@unilogcorp.com")

# Column B widened to fit the longer email addresses.
$ws.Columns.Item(2).ColumnWidth = 35.85546875

# Selection moves to D4.
$ws.Range("D4").Select()
